$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 18.59688758850098
$ws.Range("D2").Value = 176

$ws.Range("C3").Value = 17.88806915283203
$ws.Range("D3").Value = 174

$ws.Range("C4").Value = 17.99893379211426
$ws.Range("D4").Value = 176

$ws.Range("C5").Value = 17.9450511932373
$ws.Range("D5").Value = 176

$ws.Range("C6").Value = 18.16797256469727
$ws.Range("D6").Value = 177
